$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp in the title cell (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 16:10"

# --- Update per-country COVID figures (columns: B Casos totales, C Nuevos casos,
#     D Casos activos, E Recuperados, F Casos criticos, G Muertes hoy, H Muertes) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1795665
$ws.Range("C4").Value = 2135
$ws.Range("D4").Value = 519709
$ws.Range("E4").Value = 1171374
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 104582

# Alemania (row 11)
$ws.Range("B11").Value = 183113
$ws.Range("C11").Value = 94
$ws.Range("E11").Value = 9615

# Azerbaiyan (row 70)
$ws.Range("B70").Value = 5246
$ws.Range("C70").Value = 257
$ws.Range("D70").Value = 3327
$ws.Range("E70").Value = 1858
$ws.Range("G70").Value = 3
$ws.Range("H70").Value = 61

# Swap Fiyi / Curazao ordering (rows 198-199) - names swap places, figures follow
$ws.Range("A198").Value = "Curazao"
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 1

$ws.Range("A199").Value = "Fiyi"
$ws.Range("D199").Value = 15
$ws.Range("H199").Value = 0

# Swap Seychelles / Montserrat ordering (rows 210-211)
$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Swap Papua Nueva Guinea / Islas Virgenes Britanicas ordering (rows 213-214)
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
